$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7: "Experimental" property value was empty, now set to the literal text "false".
# A plain Value/Formula assignment of "false" gets auto-coerced to a Boolean cell by
# Excel, so instead write it as a formula that evaluates to the text string, then
# convert the formula to its resulting value in place - this keeps the cell a normal
# text (shared-string) cell and preserves the existing cell style.
$cell = $ws.Range("B7")
$cell.Formula = '="false"'
$cell.Copy()
$cell.PasteSpecial(-4163)

# Row 8: "Date" property value updated to new timestamp
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
